$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fix: "unnamed: 1_level_1" -> "total" ---
$ws.Range("B2").Value = "total"

# --- Body rows: drop the "situação do domicílio" / "grandes regiões"
#     sub-header-only rows, shifting every category up by one row and
#     filling in the correct figures for each. ---
$rows = @(
    @{ Row=5;  Label="urbana";       Vals=@(2.33, 6.48, 2.96, 5.32, 5.71, 4.12, 8.76, 7.6) }
    @{ Row=6;  Label="rural";        Vals=@(7.28, 13.01, 9.5, 15.48, 24.16, 15.56, 59.03, 36.76) }
    @{ Row=7;  Label="norte";        Vals=@(5.08, 16.61, 6.13, 12.15, 14.08, 9.07, 19.78, 25.5) }
    @{ Row=8;  Label="nordeste";     Vals=@(4.1, 9.130000000000001, 4.83, 8.76, 8.67, 7.66, 16.67, 14.58) }
    @{ Row=9;  Label="sudeste";      Vals=@(3.93, 12.46, 5.5, 8.58, 10.29, 7.05, 13.7, 11.64) }
    @{ Row=10; Label="sul";          Vals=@(4.82, 15.78, 6.41, 12.33, 14.17, 9.039999999999999, 20.57, 17.95) }
    @{ Row=11; Label="centro-oeste"; Vals=@(6.2, 16.43, 8.09, 13.73, 15.58, 11.23, 22.77, 20.15) }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Label
    $col = 2
    foreach ($v in $r.Vals) {
        $ws.Cells.Item($r.Row, $col).Value = $v
        $col++
    }
}

# The data that used to live in rows 12/13 has now been folded into
# rows 10/11 above, so drop the now-redundant trailing rows.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()

Write-Host "edit applied"
